$newValues = @(
    "3+2=",
    "48-8=",
    "79+6=",
    "60+24=",
    "77-76=",
    "28+17=",
    "25+50=",
    "56-56=",
    "32+38=",
    "28+10=",
    "78+11=",
    "67-54=",
    "9+72=",
    "62-48=",
    "2+5=",
    "72-42=",
    "38+13=",
    "90+6=",
    "14+49=",
    "94-59=",
    "71-8=",
    "22+30=",
    "27+69=",
    "9-8=",
    "52-40=",
    "32-9=",
    "85-14=",
    "48+36=",
    "99-46=",
    "14+38=",
    "62+7=",
    "55-24=",
    "32-5=",
    "87-6=",
    "53-2=",
    "64-16=",
    "1+71=",
    "37+55=",
    "34-20=",
    "32+58=",
    "36+8=",
    "81-71=",
    "77+4=",
    "14+11=",
    "86-4=",
    "49-39=",
    "61+26=",
    "73+19=",
    "5+23=",
    "48-35=",
    "82-12=",
    "57-47=",
    "37-31=",
    "18+3=",
    "62+32=",
    "23-17=",
    "66-56=",
    "95-8=",
    "80-11=",
    "84-34=",
    "6+88=",
    "92-87=",
    "90-72=",
    "42-15=",
    "18+38=",
    "80-73=",
    "55-11=",
    "24+58=",
    "50-35=",
    "38-18=",
    "88-34=",
    "10+59=",
    "62-23=",
    "55-42=",
    "36+59=",
    "38+12=",
    "46-31=",
    "30-5=",
    "40+24=",
    "56+13=",
    "3+56=",
    "96-81=",
    "5+75=",
    "43-9=",
    "94-15=",
    "0+37=",
    "42-38=",
    "23+14=",
    "94+2=",
    "22+40=",
    "57+22=",
    "67+2=",
    "68+13=",
    "71-55=",
    "11+28=",
    "2+90=",
    "62-28=",
    "50+38=",
    "66+2=",
    "44+1="
)

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$numRows = $t.Rows.Count
$numCols = $t.Columns.Count

$idx = 0
for ($r = 1; $r -le $numRows; $r++) {
    for ($c = 1; $c -le $numCols; $c++) {
        $cell = $t.Cell($r, $c)
        $newText = $newValues[$idx]
        $cell.Range.Text = $newText
        $idx = $idx + 1
    }
}

Write-Host "Updated $idx cells."
